$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2: Bitcoin
Set-TextValue "D2" '26.592.12'
Set-TextValue "E2" '  +4.11%  '

# Row 3: Ethereum
Set-TextValue "D3" '1.743.02'
Set-TextValue "E3" '  +4.34%  '

# Row 4: TetherUSD
Set-TextValue "D4" '0.9998'
Set-TextValue "E4" '  +0.02%  '

# Row 5: BNB
Set-TextValue "D5" '247.13'
Set-TextValue "E5" '  +3.62%  '

# Row 6: USDC
Set-TextValue "E6" '  +0.02%  '

# Row 7: XRP
Set-TextValue "D7" '0.4807'
Set-TextValue "E7" '  +0.77%  '

# Row 8: Cardano
Set-TextValue "D8" '0.2689'
Set-TextValue "E8" '  +2.83%  '

# Row 9: Dogecoin
Set-TextValue "D9" '0.06257'
Set-TextValue "E9" '  +1.39%  '

# Row 10: WrappedEther
Set-TextValue "D10" '1.742.02'
Set-TextValue "E10" '  +4.19%  '

# Row 11: TRON
Set-TextValue "D11" '0.07119'
Set-TextValue "E11" '  +2.10%  '

# Row 12: Solana
Set-TextValue "D12" '15.82'
Set-TextValue "E12" '  +6.57%  '

# Row 13: Polygon
Set-TextValue "D13" '0.6219'
Set-TextValue "E13" '  +5.59%  '

# Row 14: Polkadot
Set-TextValue "D14" '4.505'
Set-TextValue "E14" '  +2.89%  '

# Row 15: Litecoin
Set-TextValue "D15" '77.44'
Set-TextValue "E15" '  +2.80%  '

# Row 16: Dai
Set-TextValue "E16" '  +0.05%  '

# Row 17: WrappedBTC
Set-TextValue "D17" '26.590.58'
Set-TextValue "E17" '  +4.13%  '

# Row 18: BinanceUSD
Set-TextValue "D18" '1.000'
Set-TextValue "E18" '  -0.02%  '

# Row 19: ShibaInu
Set-TextValue "D19" '0.000006894'
Set-TextValue "E19" '  +2.16%  '

# Row 20: Avalanche
Set-TextValue "D20" '11.71'
Set-TextValue "E20" '  +2.32%  '

# Row 21: WrappedliquidstakedEther2.0
Set-TextValue "D21" '1.966.56'
Set-TextValue "E21" '  +4.17%  '

# Row 22: Uniswap
Set-TextValue "E22" '  +4.05%  '

# Row 23: Cosmos
Set-TextValue "D23" '8.832'
Set-TextValue "E23" '  +0.27%  '

# Row 24: Chainlink
Set-TextValue "D24" '5.343'
Set-TextValue "E24" '  +1.42%  '

# Row 25: Monero
Set-TextValue "E25" '  -0.61%  '

# Row 26: EthereumClassic
Set-TextValue "D26" '15.39'
Set-TextValue "E26" '  +2.30%  '

# Row 27: LidoDAOToken
Set-TextValue "D27" '1.819'
Set-TextValue "E27" '  +5.36%  '

# Row 28: Toncoin
Set-TextValue "D28" '1.432'
Set-TextValue "E28" '  +2.95%  '

# Row 29: BitcoinCash
Set-TextValue "D29" '107.06'
Set-TextValue "E29" '  +2.31%  '

# Row 30: InternetComputer(DFINITY)
Set-TextValue "D30" '4.016'
Set-TextValue "E30" '  +0.86%  '

# Row 31: Filecoin
Set-TextValue "D31" '3.756'
Set-TextValue "E31" '  +3.47%  '

# Row 32: Stellar
Set-TextValue "D32" '0.07864'
Set-TextValue "E32" '  -0.09%  '

# Row 33: Hedera
Set-TextValue "D33" '0.04605'
Set-TextValue "E33" '  +7.36%  '

# Row 34: Frax
Set-TextValue "B34" 'Frax'
Set-TextValue "C34" 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue "D34" '1.000'
Set-TextValue "E34" '  +0.09%  '

# Row 35: HuobiToken
Set-TextValue "B35" 'HuobiToken'
Set-TextValue "C35" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D35" '2.619'
Set-TextValue "E35" '  -0.04%  '

# Row 36: ImmutableX
Set-TextValue "B36" 'ImmutableX'
Set-TextValue "C36" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D36" '0.6446'
Set-TextValue "E36" '  +6.06%  '

# Row 37: ARBITRUM
Set-TextValue "B37" 'ARBITRUM'
Set-TextValue "C37" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D37" '0.9969'
Set-TextValue "E37" '  +4.43%  '

# Row 38: TrustWalletToken
Set-TextValue "B38" 'TrustWalletToken'
Set-TextValue "C38" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D38" '0.9469'
Set-TextValue "E38" '  +5.75%  '

# Row 39: Quant
Set-TextValue "B39" 'Quant'
Set-TextValue "C39" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D39" '113.28'
Set-TextValue "E39" '  +17.62%  '

# Row 40: RenderToken
Set-TextValue "B40" 'RenderToken'
Set-TextValue "C40" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D40" '1.994'
Set-TextValue "E40" '  +7.28%  '

# Row 41: MXToken
Set-TextValue "B41" 'MXToken'
Set-TextValue "C41" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D41" '2.422'
Set-TextValue "E41" '  -6.49%  '

# Row 42: PaxDollar
Set-TextValue "B42" 'PaxDollar'
Set-TextValue "C42" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D42" '1.004'
Set-TextValue "E42" '  +0.41%  '

# Row 43: FraxShare
Set-TextValue "B43" 'FraxShare'
Set-TextValue "C43" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" '5.764'
Set-TextValue "E43" '  +17.40%  '

# Row 44: VeChain
Set-TextValue "B44" 'VeChain'
Set-TextValue "C44" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D44" '0.01507'
Set-TextValue "E44" '  +1.89%  '

# Row 45: TheSandbox
Set-TextValue "B45" 'TheSandbox'
Set-TextValue "C45" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D45" '0.3914'
Set-TextValue "E45" '  +4.28%  '

# Row 46: Algorand
Set-TextValue "B46" 'Algorand'
Set-TextValue "C46" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D46" '0.1209'
Set-TextValue "E46" '  +7.98%  '

# Row 47: Aptos
Set-TextValue "B47" 'Aptos'
Set-TextValue "C47" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D47" '6.694'
Set-TextValue "E47" '  +7.75%  '

# Row 48: Cronos
Set-TextValue "B48" 'Cronos'
Set-TextValue "C48" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D48" '0.05328'
Set-TextValue "E48" '  +1.15%  '

# Row 49: EnergySwap
Set-TextValue "B49" 'EnergySwap'
Set-TextValue "C49" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D49" '7.975'
Set-TextValue "E49" '  +6.86%  '

# Row 50: Elrond
Set-TextValue "B50" 'Elrond'
Set-TextValue "C50" 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue "D50" '30.74'
Set-TextValue "E50" '  +2.63%  '

# Row 51: NEARProtocol
Set-TextValue "B51" 'NEARProtocol'
Set-TextValue "C51" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D51" '1.270'
Set-TextValue "E51" '  +5.50%  '
